# Apply the BOM update: a new "Test pad" (TP2) line item is inserted as row 36,
# pushing the existing KF33D / TJA1044 / USB_MICROB / 3M Pak rows down by one
# (old rows 36-39 become new rows 37-40).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 36; this shifts the old rows 36-39 (and their
# per-row formatting, such as row 39's ht="64.5") down to rows 37-40.
$ws.Rows("36:36").Insert()

# Populate the new row 36 with the Test pad BOM entry.
$ws.Range("A36").Value = "TPB2,54"
$ws.Range("B36").Value = "<b>Test pad</b>"
$ws.Range("C36").Value = "TP2"
$ws.Range("D36").Value = "B2,54"
$ws.Range("E36").Value = "TPB2,54"
$ws.Range("F36").Value = 1

# The freshly inserted row doesn't carry the table's normal cell formatting
# (borders/font used throughout the BOM), so copy the formatting (only) from
# the row right below it (the old row 36 data, now at row 37) back onto row 36.
$ws.Range("A37:F37").Copy()
$ws.Range("A36:F36").PasteSpecial(-4122)
